$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new "reset password via email" API row (row 8) ---
# (Populate the new row first so the shared-string table grows in the
# same order as the authored edit before we touch the existing H5 cell.)
$ws.Range("B8").Value = "Accounts"
$ws.Range("C8").Value = "accountsService/restMyPasswordWithEmailLink"
$ws.Range("D8").Value = "to reset password via email"
$ws.Range("E8").Value = "POST"
$ws.Range("F8").Value = "{`n ""email"":""zeeshanahmedd0010@gmail.com""`n}"
$ws.Range("G8").Value = "TEXT"
$ws.Range("H8").Value = "this send an email for password varaification"

# Post-format cell (F8) wraps text like the other post-format cells (F5/F7)
$ws.Range("F8").WrapText = $true

# Service-name cell (C8) gets the small grey Arial font used for the new entry
$ws.Range("C8").Font.Name = "Arial"
$ws.Range("C8").Font.Size = 9
$ws.Range("C8").Font.Color = 5263440

# Row height for the new row
$ws.Range("B8").EntireRow.RowHeight = 45

# --- Update the "create account" row (row 5) so it explains the new verification flow ---
$ws.Range("H5").Value = "this will create account in firebase and will send the email varification link and will not allow until u verifiy your account"

# --- Move the active selection like the authored workbook ---
$ws.Range("B9").Select()
